# Scheduled-runner update: refresh cached market-price figures (columns
# H-N: currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) across the Anima_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1015.2432
$ws.Range("I6").Value = 95.31035
$ws.Range("K6").Value = 285.93105
$ws.Range("M6").Value = -173.93105
$ws.Range("H8").Value = 255.71428
$ws.Range("I8").Value = 255.71428
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 767.14284
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -628.14284
$ws.Range("N8").ClearContents()
$ws.Range("H33").Value = 2395.6667
$ws.Range("I33").Value = 1445.125
$ws.Range("K33").Value = 1445.125
$ws.Range("M33").Value = -1216.125
$ws.Range("H112").Value = 5243.5957
$ws.Range("J112").Value = 5710.2095
$ws.Range("L112").Value = 17130.6285
$ws.Range("N112").Value = -19346.6285
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H129").Value = 964.15625
$ws.Range("J129").Value = 1113.551
$ws.Range("L129").Value = 3340.653
$ws.Range("N129").Value = -13340.653
$ws.Range("H138").Value = 2349.037
$ws.Range("I138").Value = 1921
$ws.Range("J138").Value = 2666.6128
$ws.Range("K138").Value = 5763
$ws.Range("L138").Value = 7999.8384
$ws.Range("M138").Value = -623
$ws.Range("N138").Value = -18279.8384
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4464.0713
$ws.Range("I45").Value = 4512.4443
$ws.Range("J45").Value = 4377
$ws.Range("K45").Value = 4512.4443
$ws.Range("L45").Value = 4377
$ws.Range("M45").Value = -4135.4443
$ws.Range("N45").Value = -5131
$ws.Range("H102").Value = 3992.8572
$ws.Range("I102").Value = 4341.6665
$ws.Range("J102").Value = 1900
$ws.Range("K102").Value = 4341.6665
$ws.Range("L102").Value = 1900
$ws.Range("M102").Value = -2719.6665
$ws.Range("N102").Value = -5144
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 19390
$ws.Range("J74").Value = 19390
$ws.Range("L74").Value = 19390
$ws.Range("N74").Value = -21262
$ws.Range("H77").Value = 19390
$ws.Range("J77").Value = 19390
$ws.Range("L77").Value = 58170
$ws.Range("N77").Value = -67530
$ws.Range("H105").Value = 2947.2632
$ws.Range("I105").Value = 2941.0588
$ws.Range("K105").Value = 2941.0588
$ws.Range("M105").Value = -1194.0588
$ws.Range("H134").Value = 2193.838
$ws.Range("I134").Value = 1958.0869
$ws.Range("J134").Value = 2581.1428
$ws.Range("K134").Value = 5874.2607
$ws.Range("L134").Value = 7743.428400000001
$ws.Range("M134").Value = -3339.2607
$ws.Range("N134").Value = -12813.4284
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2048.4348
$ws.Range("I99").Value = 2266.6667
$ws.Range("J99").Value = 2015.7
$ws.Range("K99").Value = 2266.6667
$ws.Range("L99").Value = 2015.7
$ws.Range("M99").Value = -768.6667000000002
$ws.Range("N99").Value = -5011.7
$ws.Range("H126").Value = 2048.4348
$ws.Range("I126").Value = 2266.6667
$ws.Range("J126").Value = 2015.7
$ws.Range("K126").Value = 6800.000100000001
$ws.Range("L126").Value = 6047.1
$ws.Range("M126").Value = -4330.000100000001
$ws.Range("N126").Value = -10987.1
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1332.2222
$ws.Range("J31").Value = 1123.75
$ws.Range("L31").Value = 3371.25
$ws.Range("N31").Value = -3947.25
$ws.Range("H44").Value = 261.625
$ws.Range("I44").Value = 261.625
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 784.875
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -386.875
$ws.Range("N44").ClearContents()
$ws.Range("H51").Value = 1165.2
$ws.Range("I51").Value = 768
$ws.Range("J51").Value = 1235.2941
$ws.Range("K51").Value = 2304
$ws.Range("L51").Value = 3705.8823
$ws.Range("M51").Value = -1844
$ws.Range("N51").Value = -4625.8823
$ws.Range("H118").Value = 2873.2559
$ws.Range("I118").Value = 1322.5
$ws.Range("J118").Value = 3032.3076
$ws.Range("K118").Value = 3967.5
$ws.Range("L118").Value = 9096.9228
$ws.Range("M118").Value = -2724.5
$ws.Range("N118").Value = -11582.9228
$ws.Range("H131").Value = 1134.7407
$ws.Range("J131").Value = 1187.52
$ws.Range("L131").Value = 3562.56
$ws.Range("N131").Value = -13642.56
$ws.Range("H134").Value = 7079.1113
$ws.Range("I134").Value = 3961.4
$ws.Range("J134").Value = 7787.6816
$ws.Range("K134").Value = 11884.2
$ws.Range("L134").Value = 23363.0448
$ws.Range("M134").Value = -6814.200000000001
$ws.Range("N134").Value = -33503.0448
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 55159.332
$ws.Range("J127").Value = 55159.332
$ws.Range("L127").Value = 55159.332
$ws.Range("N127").Value = -65079.332
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2310.8572
$ws.Range("I7").Value = 2310.8572
$ws.Range("K7").Value = 2310.8572
$ws.Range("M7").Value = -2198.8572
$ws.Range("H55").Value = 842.38464
$ws.Range("I55").Value = 696.375
$ws.Range("K55").Value = 696.375
$ws.Range("M55").Value = -523.375
$ws.Range("H122").Value = 3548.4
$ws.Range("I122").Value = 3182.0908
$ws.Range("K122").Value = 9546.2724
$ws.Range("M122").Value = -7096.2724
$ws.Range("H126").Value = 2310.8572
$ws.Range("I126").Value = 2310.8572
$ws.Range("K126").Value = 6932.571599999999
$ws.Range("M126").Value = -4462.571599999999
$ws.Range("H132").Value = 2276.3428
$ws.Range("I132").Value = 1533.1818
$ws.Range("J132").Value = 3534
$ws.Range("K132").Value = 4599.5454
$ws.Range("L132").Value = 10602
$ws.Range("M132").Value = -2069.5454
$ws.Range("N132").Value = -15662
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 18733.75
$ws.Range("J101").Value = 18733.75
$ws.Range("L101").Value = 18733.75
$ws.Range("N101").Value = -25223.75
$ws.Range("H132").Value = 932.6795
$ws.Range("I132").Value = 662.38336
$ws.Range("J132").Value = 1833.6666
$ws.Range("K132").Value = 1987.15008
$ws.Range("L132").Value = 5500.9998
$ws.Range("M132").Value = 542.8499199999999
$ws.Range("N132").Value = -10560.9998
